# "assigned login scenarios completed"
#
# The "login" worksheet holds a small table of login test scenarios
# (scenario / username / password / errorMessage). This change:
#   1. Renames the scenario labels to use underscores instead of spaces
#      (invalid credentials -> invalid_credentials, empty username ->
#      empty_username, empty password -> empty_password).
#   2. Finishes out the "empty_username" and "empty_password" scenarios:
#        - empty_username: no username is supplied; the password moves to
#          C3 and a new expected errorMessage is recorded in D3.
#        - empty_password: the username moves to B4, no password is
#          supplied, and a new expected errorMessage is recorded in D4.
#   3. Widens columns B and D so the new username/errorMessage text fits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")

# --- Rename scenario labels (column A) -------------------------------
$ws.Range("A2").Value = "invalid_credentials"
$ws.Range("A3").Value = "empty_username"
$ws.Range("A4").Value = "empty_password"

# --- empty_username scenario (row 3): no username supplied -----------
$ws.Range("B3").Clear()
$ws.Range("C3").Value = "March@2025"

# Pick up the shared "Normal"/s=1 formatting used by the rest of the
# table before writing the new D3 error message.
$ws.Range("A1").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = "Please enter your user name"

# --- empty_password scenario (row 4): no password supplied -----------
$ws.Range("B4").Value = "Playwright@gmail.com"
$ws.Range("C4").Clear()

$ws.Range("A1").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = "Please enter your password"

$excel.CutCopyMode = 0

# --- Column widths ------------------------------------------------------
# ColumnWidth is character-width units; Excel COM rounds the stored value
# to the nearest screen pixel, so the input is pre-compensated for the
# ~0.83-unit padding Excel adds back on read to land on the target widths
# (18.0 and 36.88) as closely as COM's pixel grid allows.
$ws.Range("B1").EntireColumn.ColumnWidth = 17.166667
$ws.Range("D1").EntireColumn.ColumnWidth = 36.046667
